# homework 2 results updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Total Vize" column K needs its own width (matches 17 "characters" once
# stored in the OOXML col element)
$ws.Columns.Item(11).ColumnWidth = 16.1666666667

# New "Total Vize" column header in K1 - copy the existing header formatting
# (bold white-on-grey) from I1 so the new header matches the others, then set
# its text.
$ws.Range("I1").Copy()
$ws.Range("K1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K1").Value = "Total Vize"

# K2 is a standalone formula (matches the un-shared formula pattern used for I2)
$ws.Range("K2").Formula = "=(I2/40)*100"

# K3:K66 mirror the existing shared-formula block used by I3:I66
$ws.Range("K3:K66").Formula = "=(I3/40)*100"

# K67:K118 mirror the existing shared-formula block used by I67:I118
$ws.Range("K67:K118").Formula = "=(I67/40)*100"

# Participant in row 108 was a placeholder "x" -> now has a real name
$ws.Range("A108").Value = "şevval özmen"

# Row 113 received a late "Ödev 2" (F) grade of 100
$ws.Range("F113").Value = 100

# Restore the scrolled view / active selection from the edit session
$ws.Range("G107").Select()
